$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.107983
$ws.Range("H2").Value = 0.323949
$ws.Range("M2").Value = 0.232947
$ws.Range("N2").Value = 0.698841
$ws.Range("O2").Value = 0.2572219815457369
$ws.Range("P2").Value = 0.2572219815457369
$ws.Range("Q2").Value = 0.025154315901
$ws.Range("R2").Value = 0.226388843109
$ws.Range("S2").Value = 0.2572219815457369
$ws.Range("T2").Value = 0.2572219815457369

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.107983
$ws.Range("H3").Value = 0.323949
$ws.Range("M3").Value = 0.6726793333333333
$ws.Range("N3").Value = 2.018038
$ws.Range("O3").Value = 0.7427780184542632
$ws.Range("P3").Value = 0.7427780184542632
$ws.Range("Q3").Value = 0.07263793245133332
$ws.Range("R3").Value = 0.6537413920619999
$ws.Range("S3").Value = 0.7427780184542632
$ws.Range("T3").Value = 0.7427780184542632
